# mrange_analysis.xlsx edit:
#  - add a "PYTHON" header in K1 (formatted like the other header cells,
#    i.e. bold / yellow fill / left+right border only - same look as the
#    existing header style used for A1:J1 but without top/bottom borders)
#  - widen column K to fit the new formula text
#  - apply an AutoFilter over A1:K65 that filters column H (the "DUMP"
#    column, 0-based field index 7) down to rows whose value is "False",
#    which hides the "True" rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K1 header cell -------------------------------------------------
$ws.Range("K1").Value = "PYTHON"

# Match the header formatting (bold font, yellow fill) used by A1:J1,
# but with only left/right borders (no top/bottom) - copy the existing
# header format then strip the top/bottom border edges.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Borders.Item(8).LineStyle = -4142
$ws.Range("K1").Borders.Item(9).LineStyle = -4142
$excel.CutCopyMode = $false

# --- column width -----------------------------------------------------
$ws.Columns("K").ColumnWidth = 55

# --- AutoFilter ---------------------------------------------------------
# Field 8 = column H (1-based within the A1:K65 range), keep rows equal to
# "False" -> hides every row whose H value is "True".
$rng = $ws.Range("A1:K65")
$rng.AutoFilter(8, @("False"))
